# The "Priority List" document currently lists three consecutive
# list-paragraphs in this order:
#   "Course Feedback."
#   "General Materials."
#   "Course Materials."
# The commit reorders them to:
#   "General Materials."
#   "Course Materials."
#   "Course Feedback."
#
# All three paragraphs share identical formatting (ListParagraph style,
# the same numbering definition, and a single plain run each), so the
# reorder can be achieved simply by swapping the text content of the
# runs in place rather than moving/duplicating paragraphs.

$d = $word.ActiveDocument

function Clean-Text($s) {
    return $s.TrimEnd([char]13, [char]7)
}

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ((Clean-Text $p.Range.Text) -eq "Course Feedback.") {
        $target = $i
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the 'Course Feedback.' paragraph"
}

$p1 = $d.Paragraphs.Item($target)
$p2 = $d.Paragraphs.Item($target + 1)
$p3 = $d.Paragraphs.Item($target + 2)

$p1Text = Clean-Text $p1.Range.Text
$p2Text = Clean-Text $p2.Range.Text
$p3Text = Clean-Text $p3.Range.Text

if ($p1Text -ne "Course Feedback." -or $p2Text -ne "General Materials." -or $p3Text -ne "Course Materials.") {
    throw "Unexpected paragraph sequence around 'Course Feedback.' ($p1Text / $p2Text / $p3Text)"
}

$p1.Range.Text = "General Materials."
$p2.Range.Text = "Course Materials."
$p3.Range.Text = "Course Feedback."
